$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$ws1.Range("F3").Value = 253
$ws1.Range("F4").Value = 551
$ws1.Range("F5").Value = 2422
$ws1.Range("F7").Value = 140
$ws1.Range("F9").Value = 228
$ws1.Range("F10").Value = 4865
$ws1.Range("F11").Value = 6220
$ws1.Range("F12").Value = 867
$ws1.Range("F14").Value = 1377
$ws1.Range("F15").Value = 1333
$ws1.Range("F16").Value = 530
$ws1.Range("F17").Value = 6757
$ws1.Range("F18").Value = 374
$ws1.Range("F21").Value = 4518
$ws1.Range("F22").Value = 357
$ws1.Range("F23").Value = 21
$ws1.Range("F24").Value = 745
$ws1.Range("F25").Value = 2180
$ws1.Range("F26").Value = 1220
$ws1.Range("F27").Value = 402
$ws1.Range("F29").Value = 159
$ws1.Range("F30").Value = 72
$ws1.Range("F31").Value = 57
$ws1.Range("F35").Value = 1948
$ws1.Range("F36").Value = 188
$ws1.Range("F37").Value = 469
$ws1.Range("F39").Value = 1300
$ws1.Range("F42").Value = 58
$ws1.Range("F44").Value = 1325
$ws1.Range("F49").Value = 46

# Sheet 2
$ws2.Range("F7").Value = 235
$ws2.Range("F10").Value = 8
$ws2.Range("F12").Value = 352
$ws2.Range("F13").Value = 233
$ws2.Range("F15").Value = 150
$ws2.Range("F20").Value = 11
$ws2.Range("F25").Value = 297
$ws2.Range("F26").Value = 244
$ws2.Range("F27").Value = 14

# Sheet 3
$ws3.Range("F6").Value = 1608
$ws3.Range("F7").Value = 514
$ws3.Range("F9").Value = 1185
$ws3.Range("F10").Value = 1194
$ws3.Range("F11").Value = 1677
$ws3.Range("F12").Value = 2006
$ws3.Range("F13").Value = 470
$ws3.Range("F14").Value = 371

# Sheet 4
$ws4.Range("F2").Value = 1608
$ws4.Range("F3").Value = 551
$ws4.Range("F4").Value = 514
$ws4.Range("F5").Value = 2423
$ws4.Range("F6").Value = 1185
$ws4.Range("F7").Value = 1677
$ws4.Range("F8").Value = 228
$ws4.Range("F9").Value = 2006
$ws4.Range("F10").Value = 4865
$ws4.Range("F11").Value = 470
$ws4.Range("F13").Value = 235
$ws4.Range("F14").Value = 867
$ws4.Range("F17").Value = 1377
$ws4.Range("F18").Value = 1333
$ws4.Range("F19").Value = 530
$ws4.Range("F20").Value = 6757
$ws4.Range("F21").Value = 374
$ws4.Range("F22").Value = 371
$ws4.Range("F23").Value = 8
$ws4.Range("F25").Value = 4518
$ws4.Range("F26").Value = 357
$ws4.Range("F27").Value = 745
$ws4.Range("F28").Value = 2180
$ws4.Range("F29").Value = 1220
$ws4.Range("F30").Value = 402
$ws4.Range("F32").Value = 159
$ws4.Range("F33").Value = 57
$ws4.Range("F34").Value = 233
$ws4.Range("F38").Value = 1948
$ws4.Range("F39").Value = 188
$ws4.Range("F40").Value = 469
$ws4.Range("F42").Value = 11
$ws4.Range("F43").Value = 1300
$ws4.Range("F46").Value = 14
$ws4.Range("F48").Value = 1325
